# Updated the request template to include compliance flags.
#
# Target: slide layout "3_Title and Content" (the 4th custom layout on
# the slide master) - the layout that hosts the Project Request
# template's status / compliance placeholders.

$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(4)
$shapes = $layout.Shapes

# Resolve the three placeholders we need by their stable shape Id
# (idx="24" -> id 39, idx="29" -> id 28, idx="30" -> id 30) instead of a
# hard-coded collection position, so the script keeps working even if
# shapes get reordered.
$statusShape = $null
$complianceTop = $null
$complianceBottom = $null
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Id -eq 39) { $statusShape = $sh }
    elseif ($sh.Id -eq 28) { $complianceTop = $sh }
    elseif ($sh.Id -eq 30) { $complianceBottom = $sh }
}

# --- Placeholder idx="24" ("Status" label area) ----------------------------
# Nudge it down slightly and give it real prompt text ("This is a status")
# instead of the inherited "Click to edit Master text styles" text.
$statusShape.Top = 108.86283564566929
$statusShape.TextFrame.TextRange.Text = "This is a status"

# --- Placeholder idx="29" / idx="30" ("Compliance" flags) -----------------
# Clear the custom "Compliance" prompt text back to an empty placeholder.
$complianceTop.TextFrame.TextRange.Text = ""
$complianceBottom.TextFrame.TextRange.Text = ""
